# Apply updated Gungnir profit figures across the Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1710545.1
$ws.Range("I86").Value = 47999.8
$ws.Range("J86").Value = 2349985.8
$ws.Range("K86").Value = 47999.8
$ws.Range("L86").Value = 2349985.8
$ws.Range("M86").Value = -46876.8
$ws.Range("N86").Value = -2352231.8

$ws.Range("H89").Value = 1710545.1
$ws.Range("I89").Value = 47999.8
$ws.Range("J89").Value = 2349985.8
$ws.Range("K89").Value = 239999
$ws.Range("L89").Value = 11749929
$ws.Range("M89").Value = -234383
$ws.Range("N89").Value = -11761161

$ws.Range("H100").Value = 100000
$ws.Range("I100").Value = 100000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 100000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -99459
$ws.Range("N100").ClearContents()

$ws.Range("H106").Value = 2883.76
$ws.Range("I106").Value = 2461.3845
$ws.Range("K106").Value = 2461.3845
$ws.Range("M106").Value = -1830.3845

$ws.Range("H132").Value = 7697347.5
$ws.Range("I132").Value = 10208673
$ws.Range("J132").Value = 6414.75
$ws.Range("K132").Value = 30626019
$ws.Range("L132").Value = 19244.25
$ws.Range("M132").Value = -30623489
$ws.Range("N132").Value = -24304.25

$ws.Range("H137").Value = 1184.1621
$ws.Range("I137").Value = 1072.2142
$ws.Range("J137").Value = 1532.4445
$ws.Range("K137").Value = 3216.6426
$ws.Range("L137").Value = 4597.333500000001
$ws.Range("M137").Value = -666.6425999999997
$ws.Range("N137").Value = -9697.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1175.8823
$ws.Range("I61").Value = 1137.037
$ws.Range("J61").Value = 1325.7142
$ws.Range("K61").Value = 1137.037
$ws.Range("L61").Value = 1325.7142
$ws.Range("M61").Value = -925.037
$ws.Range("N61").Value = -1749.7142

$ws.Range("H74").Value = 1116.7742
$ws.Range("I74").Value = 1157.5862
$ws.Range("J74").Value = 525
$ws.Range("K74").Value = 1157.5862
$ws.Range("L74").Value = 525
$ws.Range("M74").Value = -283.5862
$ws.Range("N74").Value = -2273

$ws.Range("H77").Value = 1116.7742
$ws.Range("I77").Value = 1157.5862
$ws.Range("J77").Value = 525
$ws.Range("K77").Value = 5787.931
$ws.Range("L77").Value = 2625
$ws.Range("M77").Value = -1419.931
$ws.Range("N77").Value = -11361

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 1549576
$ws.Range("I132").Value = 1138.52
$ws.Range("J132").Value = 4527340.5
$ws.Range("K132").Value = 3415.56
$ws.Range("L132").Value = 13582021.5
$ws.Range("M132").Value = -885.5599999999999
$ws.Range("N132").Value = -13587081.5

$ws.Range("H136").Value = 1175.8823
$ws.Range("I136").Value = 1137.037
$ws.Range("J136").Value = 1325.7142
$ws.Range("K136").Value = 3411.111
$ws.Range("L136").Value = 3977.1426
$ws.Range("M136").Value = -861.1109999999999
$ws.Range("N136").Value = -9077.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2936.1428
$ws.Range("I134").Value = 762.4545000000001
$ws.Range("J134").Value = 10906.333
$ws.Range("K134").Value = 2287.3635
$ws.Range("L134").Value = 32718.999
$ws.Range("M134").Value = 247.6364999999996
$ws.Range("N134").Value = -37788.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1359.0577
$ws.Range("I31").Value = 965.75
$ws.Range("J31").Value = 1604.875
$ws.Range("K31").Value = 965.75
$ws.Range("L31").Value = 1604.875
$ws.Range("M31").Value = -670.75
$ws.Range("N31").Value = -2194.875

$ws.Range("H34").Value = 1359.0577
$ws.Range("I34").Value = 965.75
$ws.Range("J34").Value = 1604.875
$ws.Range("K34").Value = 965.75
$ws.Range("L34").Value = 1604.875
$ws.Range("M34").Value = -763.75
$ws.Range("N34").Value = -2008.875

$ws.Range("H58").Value = 15625869
$ws.Range("I58").Value = 22728180
$ws.Range("J58").Value = 785.8
$ws.Range("K58").Value = 22728180
$ws.Range("L58").Value = 785.8
$ws.Range("M58").Value = -22727977
$ws.Range("N58").Value = -1191.8

$ws.Range("H105").Value = 40400
$ws.Range("I105").Value = 48005
$ws.Range("J105").Value = 9980
$ws.Range("K105").Value = 48005
$ws.Range("L105").Value = 9980
$ws.Range("M105").Value = -46258
$ws.Range("N105").Value = -13474

$ws.Range("H132").Value = 8773423
$ws.Range("I132").Value = 1236.238
$ws.Range("J132").Value = 19609654
$ws.Range("K132").Value = 3708.714
$ws.Range("L132").Value = 58828962
$ws.Range("M132").Value = -1178.714
$ws.Range("N132").Value = -58834022

$ws.Range("H134").Value = 1161.3214
$ws.Range("I134").Value = 1228.7142
$ws.Range("K134").Value = 3686.1426
$ws.Range("M134").Value = -1151.1426

$ws.Range("H136").Value = 15625869
$ws.Range("I136").Value = 22728180
$ws.Range("J136").Value = 785.8
$ws.Range("K136").Value = 68184540
$ws.Range("L136").Value = 2357.4
$ws.Range("M136").Value = -68181990
$ws.Range("N136").Value = -7457.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 543.6842
$ws.Range("I40").Value = 205
$ws.Range("J40").Value = 700
$ws.Range("K40").Value = 820
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -751
$ws.Range("N40").Value = -2938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15501.75
$ws.Range("I132").Value = 12201.777
$ws.Range("J132").Value = 19744.572
$ws.Range("K132").Value = 36605.331
$ws.Range("L132").Value = 59233.716
$ws.Range("M132").Value = -34075.331
$ws.Range("N132").Value = -64293.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13702988
$ws.Range("I132").Value = 22728544
$ws.Range("J132").Value = 9040.069
$ws.Range("K132").Value = 68185632
$ws.Range("L132").Value = 27120.207
$ws.Range("M132").Value = -68183102
$ws.Range("N132").Value = -32180.207

$ws.Range("H136").Value = 29306336
$ws.Range("I136").Value = 4466681.5
$ws.Range("J136").Value = 142859040
$ws.Range("K136").Value = 13400044.5
$ws.Range("L136").Value = 428577120
$ws.Range("M136").Value = -13397494.5
$ws.Range("N136").Value = -428582220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14912.22
$ws.Range("I132").Value = 18256.084
$ws.Range("K132").Value = 54768.25199999999
$ws.Range("M132").Value = -52238.25199999999

$ws.Range("H136").Value = 1906.6666
$ws.Range("I136").Value = 982.6087
$ws.Range("J136").Value = 4942.857
$ws.Range("K136").Value = 2947.8261
$ws.Range("L136").Value = 14828.571
$ws.Range("M136").Value = -397.8261000000002
$ws.Range("N136").Value = -19928.571

